$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sponsored`nLodha Bhandup New Project - 2 & 3 BHK @ ₹ 2.29 Cr All Inc`nnewproject-bhandup.co.in`nhttps://www.newproject-bhandup.co.in › lodha › new-launch`nLodha LBS New Launch Luxurious 2 & 3 BHK apartments Starting Price from ₹ 2.29 Cr All inc. Lodha Bhandup - Prelaunch Project...`nPrice List`nGet Here Price List & Floor Plan 2, 3 BHK Starts From ₹ 2.29 Cr`nBrochure`nGet Here Type, Size, Area & Units Brochure Free"
$ws.Range("B2").Value = "AddressofChoice Realty Pvt. Ltd"
$ws.Range("C2").Value = "India"
$ws.Range("A3").Value = "Sponsored`nLödha Bhandup | Luxury 2, 2.5 & 3 BHK | Avail Pre Launch Offer`nl-bhandup.com`nhttps://www.l-bhandup.com › official`nLödha on LBS Rd an upcoming residential project with the best of amenities & views. Pre..."
$ws.Range("B3").Value = "Rioga Premium Real Estate Advisory LLP"
$ws.Range("C3").Value = "India"
$ws.Range("A4").Value = "Sponsored`nNew Launch At Bhandup, Mumbai`nlódháhomz.site`nhttps://www.lódháhomz.site › site-visit › enquire-now`nNew Launch At Bhandup — New Launch Project Presents 2 & 3 BHK Apartments At Bhandup, Mumbai Starting At ₹ 2.29 Cr*"
$ws.Range("B4").Value = "DIGITAL RUBIX"
$ws.Range("C4").Value = "India"
$ws.Range("A5").Value = "Sponsored`nLodha Bhandup`nnewlaunch-property.net`nhttps://www.newlaunch-property.net › lodha_bhandup`nLodha in Bhandup West — Get Huge Discount On Booking. Schedule Your Free Site Visit & Get Complete Project Details"
$ws.Range("B5").Value = "GTF Technologies"
$ws.Range("C5").Value = "India"
$ws.Range("A6").Value = "Sponsored`nLodha Bhandup`nproptigermumbai.com`nhttps://www.proptigermumbai.com › lodha_newlaunch › bhandup_mumbai`nNew Project Launch in Bhandup — its Big Apartments with 2/3 BHK, its prime Location along the LBS Road, Starts @ ₹ 2.29Cr* The Apartments have 2/3 BHK Configurations & are Designed to Provide Ample Space & Comfort. New Launch Project. Easy Payment Plan. Excellent Connectivity. Budget Friendly.`nPricing & Floor Plan · Download Brochure · Hiranandani Projects · Platinum Group"
$ws.Range("B6").Value = "Locon Solutions Pvt. Ltd."
$ws.Range("C6").Value = "India"
$ws.Range("A7").Value = "Sponsored`nLodha Bhandup New Launch | 2/3 BHK Starts @ ₹ 2.29 Cr*`nbhandupnewlaunch.com`nhttps://www.bhandupnewlaunch.com › 2&3bhk › luxury_homes`nPre-book Lodha Bhandup at ₹1.08 Lacs* | Easy Access to Powai & R-City Mall | EOI Open Now! Modern Living at Lodha Bhandup | 10 mins to Eastern Express Hwy | Pre-book at ₹1.08 Lacs* Free Pickup & Drop."
$ws.Range("B7").Value = "Finwizz Holdings"
$ws.Range("C7").Value = "India"
$ws.Range("A8").Value = "Sponsored`nLodha Bhandup | Lodha New Launch Bhandup`npreferred-partners.co.in`nhttp://www.preferred-partners.co.in › lodha_bhandup › book_now`nLodha Bhandup Premium 2 & 3 BHK Homes Starts ₹ 2.29 Cr* On Request at Bhandup West, Mumbai. The Apartments have 2/3 BHK Configurations & are Designed to Provide Ample Space & Comfort. Limited Deals Available.`nDownload Brochure · Pricing /Costing · Book a Site Visit · Top Facilities"
$ws.Range("B8").Value = "INVESTOXPERT ADVISORS PRIVATE LIMITED"
$ws.Range("C8").Value = "India"
